# Commit: "added info for dry ice"
#
# Update the transport method for the "ACR 2022 M1" and "ACR 2022 M2"
# sample batches (worksheet rows 153-284, column G) from "Dry shipper"
# to the new, more specific option "Dry shipper or dry ice". Excel will
# automatically add the new text as a shared string the first time it is
# used and re-point every edited cell at it, leaving every other
# "Dry shipper" cell elsewhere on the sheet untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G153:G284").Value = "Dry shipper or dry ice"

# Leave the cursor where the editor ended up after making the change.
$ws.Range("I288").Select()
